$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.756.10"
$ws.Range("E2").Value = "  +1.07%  "

$ws.Range("D3").Value = "2.302.52"
$ws.Range("E3").Value = "  -0.84%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.26%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.52%  "

$ws.Range("D9").Value = "2.298.95"
$ws.Range("E9").Value = "  -0.88%  "

$ws.Range("E10").Value = "  +0.22%  "

$ws.Range("E11").Value = "  +0.55%  "

$ws.Range("E12").Value = "  -0.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.330"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.65%  "

$ws.Range("D15").Value = "2.713.56"
$ws.Range("E15").Value = "  -0.83%  "

$ws.Range("D16").Value = "59.769.74"
$ws.Range("E16").Value = "  +1.09%  "

$ws.Range("E17").Value = "  -0.91%  "

$ws.Range("D18").Value = "2.307.27"
$ws.Range("E18").Value = "  -0.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "309.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.77%  "

$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.84%  "

$ws.Range("E25").Value = "  -3.03%  "

$ws.Range("E26").Value = "  -0.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.69%  "

$ws.Range("E28").Value = "  +1.38%  "

$ws.Range("E29").Value = "  +1.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.84%  "

$ws.Range("D32").Value = "0.0₃0720"
$ws.Range("E32").Value = "  -2.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.57%  "

$ws.Range("E34").Value = "  -2.04%  "

$ws.Range("E35").Value = "  -6.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.36%  "

$ws.Range("E38").Value = "  -0.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "314.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0941"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.564"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0487"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.37%  "

$ws.Range("E49").Value = "  +20.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0212"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.02%  "
